$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.170.29"
$ws.Range("E2").Value = "  +2.81%  "

$ws.Range("D3").Value = "2.506.89"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'579.75"
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("D6").Value = "'151.36"
$ws.Range("E6").Value = "  +4.05%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.542"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("D9").Value = "2.507.44"
$ws.Range("E9").Value = "  +2.93%  "

$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("D12").Value = "'0.363"
$ws.Range("E12").Value = "  +3.53%  "

$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").Value = "'27.49"
$ws.Range("E14").Value = "  +2.81%  "

$ws.Range("D15").Value = "'0.0000184"
$ws.Range("E15").Value = "  +2.85%  "

$ws.Range("D16").Value = "2.953.48"
$ws.Range("E16").Value = "  +4.73%  "

$ws.Range("D17").Value = "63.942.26"
$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("D18").Value = "2.509.91"
$ws.Range("E18").Value = "  +2.89%  "

$ws.Range("D19").Value = "'11.66"
$ws.Range("E19").Value = "  +3.47%  "

$ws.Range("E20").Value = "  +7.67%  "

$ws.Range("D21").Value = "'331.63"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("D22").Value = "'4.23"
$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("D23").Value = "'2.10"
$ws.Range("E23").Value = "  +20.53%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'66.72"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").Value = "'638.91"
$ws.Range("E26").Value = "  +12.95%  "

$ws.Range("E27").Value = "  +9.03%  "

$ws.Range("E28").Value = "  -0.65%  "

# Row 29/30 swap: Fetch.AI <-> WrappedeETH
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.625.60"
$ws.Range("E29").Value = "  +2.80%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.54"
$ws.Range("E30").Value = "  +6.60%  "

$ws.Range("D31").Value = "'8.45"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").Value = "'0.144"
$ws.Range("E33").Value = "  -2.84%  "

$ws.Range("D34").Value = "'1.93"
$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("E35").Value = "  +7.59%  "

$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "'0.386"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").Value = "'5.56"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").Value = "'19.04"
$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +2.97%  "

# Row 42/43 swap: Monero <-> dogwifhat
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.79"
$ws.Range("E42").Value = "  +15.76%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'148.80"
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").Value = "'150.31"
$ws.Range("E45").Value = "  +1.29%  "

$ws.Range("D46").Value = "'3.79"
$ws.Range("E46").Value = "  +2.53%  "

$ws.Range("D47").Value = "'21.29"
$ws.Range("E47").Value = "  +4.44%  "

$ws.Range("D48").Value = "'0.0549"
$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").Value = "'0.615"
$ws.Range("E49").Value = "  +2.48%  "

$ws.Range("D50").Value = "'0.0239"
$ws.Range("E50").Value = "  +3.21%  "

$ws.Range("D51").Value = "'0.0928"
$ws.Range("E51").Value = "  -0.05%  "
